$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'29.495.25"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.21%  '
$ws.Range('D3').Value = "'1.970.80"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +3.73%  '
$ws.Range('D4').Value = "'1.003"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = "'326.57"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.42%  '
$ws.Range('E6').Value = '  +0.17%  '
$ws.Range('D7').Value = "'0.4650"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.45%  '
$ws.Range('D8').Value = "'0.3903"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').Value = "'46.21"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.11%  '
$ws.Range('D10').Value = "'0.07936"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.65%  '
$ws.Range('D11').Value = "'0.9876"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.24%  '
$ws.Range('D12').Value = "'22.81"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +4.53%  '
$ws.Range('D13').Value = "'1.944.71"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.86%  '
$ws.Range('D14').Value = "'7.160"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.38%  '
$ws.Range('D15').Value = "'5.800"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.91%  '
$ws.Range('D16').Value = "'0.07091"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.19%  '
$ws.Range('D17').Value = "'87.65"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.53%  '
$ws.Range('D18').Value = "'1.004"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.13%  '
$ws.Range('D19').Value = "'0.000009926"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.63%  '
$ws.Range('D20').Value = "'17.23"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.93%  '
$ws.Range('D21').Value = "'1.003"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.21%  '
$ws.Range('D22').Value = "'29.517.37"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.28%  '
$ws.Range('D23').Value = "'5.539"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.38%  '
$ws.Range('D24').Value = "'11.14"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.38%  '
$ws.Range('D25').Value = "'2.194.49"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.28%  '
$ws.Range('D26').Value = "'2.105"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.18%  '
$ws.Range('D27').Value = "'158.63"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.77%  '
$ws.Range('D28').Value = "'19.51"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.31%  '
$ws.Range('D29').Value = "'5.788"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.71%  '
$ws.Range('D30').Value = "'119.56"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.00%  '
$ws.Range('D31').Value = "'1.887"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.21%  '
$ws.Range('D32').Value = "'0.09415"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.99%  '
$ws.Range('D33').Value = "'0.8763"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.83%  '
$ws.Range('D34').Value = "'5.218"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.59%  '
$ws.Range('D35').Value = "'1.319"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.48%  '
$ws.Range('D36').Value = "'3.129"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.35%  '
$ws.Range('D37').Value = "'0.05802"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.60%  '
$ws.Range('D38').Value = "'1.159"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.65%  '
$ws.Range('D39').Value = "'0.02104"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.96%  '
$ws.Range('D40').Value = "'0.5706"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.15%  '
$ws.Range('D41').Value = "'7.697"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.47%  '
$ws.Range('D42').Value = "'0.1795"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.23%  '
$ws.Range('D43').Value = "'9.615"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.97%  '
$ws.Range('D44').Value = "'2.758"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +7.92%  '
$ws.Range('D45').Value = "'0.000002823"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +48.22%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = "'11.69"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.90%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').Value = "'0.5325"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.46%  '
$ws.Range('D48').Value = "'2.148"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.61%  '
$ws.Range('D49').Value = "'0.06913"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.41%  '
$ws.Range('E50').Value = '  -1.54%  '
$ws.Range('D51').Value = "'112.41"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.62%  '
